# "Added in customer File" - append the customer's order history rows
# (Date / Restaurant Name / Total Spent / Rewards Earned) below the
# existing header block on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=5;  A="3/2/2020'";  B="Flagstaff House";      C=80.34; D=2},
    @{Row=6;  A="3/6/2020'";  B="Chimera";               C=30.2;  D=1},
    @{Row=7;  A="3/8/2020'";  B="The Med";                C=15.3;  D=1},
    @{Row=8;  A="3/12/2020'"; B="Corrida";                C=46.78; D=2},
    @{Row=9;  A="3/13/2020'"; B="The Kitchen";            C=20.13; D=1},
    @{Row=10; A="3/16/2020'"; B="Santo";                  C=33.85; D=1},
    @{Row=11; A="3/19/2020'"; B="Black Cat";              C=57.99; D=2},
    @{Row=12; A="3/21/2020'"; B="Dushanbe Teahouse";      C=51.51; D=2},
    @{Row=13; A="3/22/2020'"; B="Snooze";                 C=22.22; D=1},
    @{Row=14; A="3/30/2020'"; B="Lucile's Creole Café";   C=36.98; D=1},
    @{Row=15; A="3/30/2020'"; B="The Sink";                C=20.13; D=1}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# Column B ("Restaurant Name") grew wider to fit "Lucile's Creole Café" /
# "Dushanbe Teahouse"; columns D and E ("Rewards Earned" / "Rewards
# Claimed") now hold data too and pick up a best-fit width.
$ws.Columns.Item(2).ColumnWidth = 16.833333333333332
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666

$ws.Range("F12").Select() | Out-Null
